$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.060.38'
$ws.Range('E2').Value = '  -3.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.843.88'
$ws.Range('E3').Value = '  -2.11%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7026'
$ws.Range('E5').Value = '  -5.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '237.57'
$ws.Range('E6').Value = '  -2.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9998'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3039'
$ws.Range('E8').Value = '  -3.86%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07448'
$ws.Range('E9').Value = '  +3.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.29'
$ws.Range('E10').Value = '  -6.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08109'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7236'
$ws.Range('E12').Value = '  -4.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.238'
$ws.Range('E13').Value = '  -3.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.822.34'
$ws.Range('E14').Value = '  +0.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '88.95'
$ws.Range('E15').Value = '  -3.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '28.831.36'
$ws.Range('E16').Value = '  -3.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.788'
$ws.Range('E17').Value = '  -5.73%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '239.93'
$ws.Range('E18').Value = '  -4.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007648'
$ws.Range('E19').Value = '  -2.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.01'
$ws.Range('E20').Value = '  -4.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9989'
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9999'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.033.58'
$ws.Range('E23').Value = '  -3.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.536'
$ws.Range('E24').Value = '  -5.75%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1460'
$ws.Range('E25').Value = '  -6.74%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '161.02'
$ws.Range('E26').Value = '  -2.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.938'
$ws.Range('E27').Value = '  -3.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.96'
$ws.Range('E28').Value = '  -3.97%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.934'
$ws.Range('E29').Value = '  -5.27%  '
$ws.Range('E30').Value = '  -7.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.476'
$ws.Range('E31').Value = '  -3.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.492'
$ws.Range('E32').Value = '  -2.85%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.006'
$ws.Range('E33').Value = '  -5.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05160'
$ws.Range('E34').Value = '  -4.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.182'
$ws.Range('E35').Value = '  -5.86%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7064'
$ws.Range('E36').Value = '  -6.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.026'
$ws.Range('E37').Value = '  +3.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.648'
$ws.Range('E38').Value = '  -2.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01869'
$ws.Range('E39').Value = '  -4.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.673'
$ws.Range('E40').Value = '  -3.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8992'
$ws.Range('E41').Value = '  +3.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.951'
$ws.Range('E42').Value = '  -1.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4279'
$ws.Range('E43').Value = '  -6.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.058.13'
$ws.Range('E44').Value = '  -4.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '69.88'
$ws.Range('E45').Value = '  -3.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9997'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.62'
$ws.Range('E47').Value = '  -2.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.749'
$ws.Range('E48').Value = '  -6.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.165'
$ws.Range('E49').Value = '  -4.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.031'
$ws.Range('E50').Value = '  -7.63%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.964.39'
$ws.Range('E51').Value = '  -3.66%  '
